$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L2").Value = 1.43
$ws.Range("O2").Value = 1.36

$ws.Range("M3").Value = 1.03

$ws.Range("M4").Value = 1.03

$ws.Range("J5").Value = 1.03
$ws.Range("K5").Value = 980
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 1.28
$ws.Range("R5").Value = 1.18

$ws.Range("G6").Value = 8.4
$ws.Range("H6").Value = 1.42
$ws.Range("I6").Value = 1.56
$ws.Range("K6").Value = 5.8
$ws.Range("L6").Value = 1.22
$ws.Range("Q6").Value = 1.61
$ws.Range("T6").Value = 1.84
$ws.Range("U6").Value = 1.9
$ws.Range("V6").Value = 2.78
$ws.Range("W6").Value = 1.13

$ws.Range("I7").Value = 2.24
$ws.Range("N7").Value = 5.9
$ws.Range("P7").Value = 2.66
$ws.Range("R7").Value = 1.67
$ws.Range("T7").Value = 1.49
$ws.Range("U7").Value = 2.64
$ws.Range("W7").Value = 1.38
$ws.Range("Y7").Value = 18.5
$ws.Range("Z7").Value = 21
$ws.Range("AA7").Value = 32
$ws.Range("AC7").Value = 11.5
$ws.Range("AF7").Value = 34
$ws.Range("AG7").Value = 18
$ws.Range("AI7").Value = 30
$ws.Range("AO7").Value = 11

$ws.Range("G8").Value = 2.32
$ws.Range("I8").Value = 3.75
$ws.Range("J8").Value = 3.45
$ws.Range("P8").Value = 2.42
$ws.Range("R8").Value = 1.57
$ws.Range("S8").Value = 2.4
$ws.Range("T8").Value = 1.54
$ws.Range("V8").Value = 1.38
$ws.Range("W8").Value = 1.76

$ws.Range("F9").Value = 5.5
$ws.Range("I9").Value = 1.74
$ws.Range("L9").Value = 1.3
$ws.Range("S9").Value = 2.92
$ws.Range("T9").Value = 1.81
$ws.Range("V9").Value = 2.34
$ws.Range("X9").Value = 17.5
$ws.Range("Y9").Value = 10.5
$ws.Range("AL9").Value = 100

$ws.Range("G10").Value = 12.5
$ws.Range("H10").Value = 1.26
$ws.Range("I10").Value = 1.32
$ws.Range("J10").Value = 6.8
$ws.Range("K10").Value = 8.199999999999999
$ws.Range("N10").Value = 8.6
$ws.Range("O10").Value = 1.1
$ws.Range("P10").Value = 3.55
$ws.Range("Q10").Value = 1.31
$ws.Range("R10").Value = 2.04
$ws.Range("S10").Value = 1.76
$ws.Range("T10").Value = 1.68
$ws.Range("U10").Value = 2.18
$ws.Range("V10").Value = 4
$ws.Range("W10").Value = 1.08
$ws.Range("Y10").Value = 17.5
$ws.Range("Z10").Value = 13
$ws.Range("AA10").Value = 12.5
$ws.Range("AC10").Value = 18.5
$ws.Range("AD10").Value = 13
$ws.Range("AE10").Value = 14
$ws.Range("AF10").Value = 140
$ws.Range("AJ10").Value = 380
$ws.Range("AK10").Value = 150
$ws.Range("AL10").Value = 110
$ws.Range("AM10").Value = 110
$ws.Range("AO10").Value = 3.3

$ws.Range("F11").Value = 4.3
$ws.Range("I11").Value = 1.77
$ws.Range("J11").Value = 3.7
$ws.Range("L11").Value = 1.22
$ws.Range("N11").Value = 2.08
$ws.Range("P11").Value = 2.08
$ws.Range("Q11").Value = 1.51
$ws.Range("S11").Value = 2.28
$ws.Range("X11").Value = 34
$ws.Range("Y11").Value = 16
$ws.Range("Z11").Value = 17.5
$ws.Range("AA11").Value = 25
$ws.Range("AB11").Value = 34
$ws.Range("AC11").Value = 15
$ws.Range("AD11").Value = 15
$ws.Range("AE11").Value = 24
$ws.Range("AG11").Value = 30
$ws.Range("AH11").Value = 26
$ws.Range("AI11").Value = 44

$ws.Range("G12").Value = 1.39
$ws.Range("J12").Value = 5.2
$ws.Range("K12").Value = 7.2
$ws.Range("L12").Value = 1.18
$ws.Range("P12").Value = 2.62
$ws.Range("Q12").Value = 1.43
$ws.Range("R12").Value = 1.65
$ws.Range("S12").Value = 2.22
$ws.Range("T12").Value = 1.94
$ws.Range("U12").Value = 1.86
$ws.Range("W12").Value = 3.5
$ws.Range("AC12").Value = 18
